$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking speaker_variant values stay text (not auto-converted to numbers)
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C32").NumberFormat = "@"

# Update existing rows 2-32: new B (id) / C (speaker_variant) values
$ws.Range("B2").Value = "#ren"
$ws.Range("C2").Value = "Ren"
$ws.Range("B3").Value = "#sol"
$ws.Range("C3").Value = "Sol"
$ws.Range("B4").Value = "#amb"
$ws.Range("C4").Value = "Amb"
$ws.Range("B5").Value = "#beul"
$ws.Range("C5").Value = "Beul"
$ws.Range("B6").Value = "#laf"
$ws.Range("C6").Value = "Laf"
$ws.Range("B7").Value = "#ed"
$ws.Range("C7").Value = "Ed"
$ws.Range("B8").Value = "#cier"
$ws.Range("C8").Value = "Cier"
$ws.Range("B9").Value = "#laf-uyt"
$ws.Range("C9").Value = "Laf uyt"
$ws.Range("B10").Value = "#3"
$ws.Range("C10").Value = "3"
$ws.Range("B11").Value = "#fuent"
$ws.Range("C11").Value = "Fuent"
$ws.Range("B12").Value = "#vit"
$ws.Range("C12").Value = "Vit"
$ws.Range("B13").Value = "#nem"
$ws.Range("C13").Value = "Nem"
$ws.Range("B14").Value = "#iean"
$ws.Range("C14").Value = "Iean"
$ws.Range("B15").Value = "#1"
$ws.Range("C15").Value = "1"
$ws.Range("B16").Value = "#averg.uyt"
$ws.Range("C16").Value = "Averg.uyt"
$ws.Range("B17").Value = "#fue"
$ws.Range("C17").Value = "Fue"
$ws.Range("B18").Value = "#4"
$ws.Range("C18").Value = "4"
$ws.Range("B19").Value = "#con"
$ws.Range("C19").Value = "Con"
$ws.Range("B20").Value = "#coning"
$ws.Range("C20").Value = "Coning"
$ws.Range("B21").Value = "#vitry"
$ws.Range("C21").Value = "Vitry"
$ws.Range("B22").Value = "#lafin"
$ws.Range("C22").Value = "Lafin"
$ws.Range("B23").Value = "#lafor"
$ws.Range("C23").Value = "Lafor"
$ws.Range("B24").Value = "#bir"
$ws.Range("C24").Value = "Bir"
$ws.Range("B25").Value = "#phil"
$ws.Range("C25").Value = "Phil"
$ws.Range("B26").Value = "#coningin"
$ws.Range("C26").Value = "Coningin"
$ws.Range("B27").Value = "#ier"
$ws.Range("C27").Value = "Ier"
$ws.Range("B28").Value = "#piet"
$ws.Range("C28").Value = "Piet"
$ws.Range("B29").Value = "#averg"
$ws.Range("C29").Value = "Averg"
$ws.Range("B30").Value = "#raeds"
$ws.Range("C30").Value = "Raeds"
$ws.Range("B31").Value = "#laforce"
$ws.Range("C31").Value = "Laforce"
$ws.Range("B32").Value = "#2"
$ws.Range("C32").Value = "2"

# Clear the is_prefered column for rows 2-21 (previously marked "x")
$ws.Range("D2:D21").ClearContents()

# Append new rows 33-36
$ws.Range("A33").Value = "https://www.dbnl.org/tekst/roel018biro03_01"
$ws.Range("B33").Value = "#nic"
$ws.Range("C33").Value = "Nic"
$ws.Range("A34").Value = "https://www.dbnl.org/tekst/roel018biro03_01"
$ws.Range("B34").Value = "#sold"
$ws.Range("C34").Value = "Sold"
$ws.Range("A35").Value = "https://www.dbnl.org/tekst/roel018biro03_01"
$ws.Range("B35").Value = "#con.-uyt"
$ws.Range("C35").Value = "Con. uyt"
$ws.Range("A36").Value = "https://www.dbnl.org/tekst/roel018biro03_01"
$ws.Range("B36").Value = "#lafo"
$ws.Range("C36").Value = "Lafo"
